$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.024.95'
$ws.Range("E2").Value = '  +2.69%  '
$ws.Range("D3").Value = '3.720.05'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.77'
$ws.Range("E5").Value = '  +9.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.89'
$ws.Range("E6").Value = '  +13.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.640'
$ws.Range("E7").Value = '  +4.35%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.728'
$ws.Range("E9").Value = '  +5.20%  '
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.52'
$ws.Range("E11").Value = '  +19.58%  '
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.53'
$ws.Range("E13").Value = '  +2.13%  '
$ws.Range("D14").Value = '4.320.74'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").Value = '3.728.76'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("E16").Value = '  +4.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.54'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  +3.08%  '
$ws.Range("D20").Value = '68.952.14'
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '413.20'
$ws.Range("E21").Value = '  +3.65%  '
$ws.Range("E22").Value = '  +4.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '90.28'
$ws.Range("E23").Value = '  +4.57%  '
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.48'
$ws.Range("E25").Value = '  +10.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.01'
$ws.Range("E26").Value = '  +4.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.82'
$ws.Range("E27").Value = '  +5.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.03'
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.81'
$ws.Range("E29").Value = '  +5.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.98'
$ws.Range("E30").Value = '  +2.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  +6.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.80'
$ws.Range("E32").Value = '  +4.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '650.22'
$ws.Range("E33").Value = '  +12.65%  '
$ws.Range("E34").Value = '  +8.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '46.29'
$ws.Range("E35").Value = '  +9.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '66.94'
$ws.Range("E36").Value = '  +4.80%  '
$ws.Range("E37").Value = '  -3.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.416'
$ws.Range("E38").Value = '  +7.48%  '
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  +6.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.08'
$ws.Range("E42").Value = '  +4.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0450'
$ws.Range("E43").Value = '  +5.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.63'
$ws.Range("E44").Value = '  +5.48%  '
$ws.Range("D45").Value = '2.907.06'
$ws.Range("E45").Value = '  +8.16%  '
$ws.Range("E46").Value = '  +6.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.29'
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '143.41'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.09'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.59'
$ws.Range("E51").Value = '  -7.97%  '
